$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.975;  "C2" = 0.999;  "D2" = 0.982;  "E2" = 0.9429999999999999;
    "F2" = 0.481;  "G2" = 0.281;  "H2" = 0.355;  "I2" = 0.725;
    "J2" = 0.862;  "K2" = 0.788;

    "B3" = 0.709;  "C3" = 0.681;  "D3" = 0.801;  "E3" = 0.272;
    "F3" = 0.547;  "G3" = 0.393;  "H3" = 0.458;  "I3" = 0.756;
    "J3" = 0.852;  "K3" = 0.801;

    "B4" = 0.929;  "C4" = 0.982;  "D4" = 0.949;  "E4" = 0.833;
    "F4" = 0.541;  "G4" = 0.371;  "H4" = 0.44;   "I4" = 0.75;
    "J4" = 0.857;  "K4" = 0.8;

    "B5" = 0.975;  "C5" = 0.997;  "D5" = 0.982;  "E5" = 0.945;
    "F5" = 0.507;  "G5" = 0.382;  "H5" = 0.436;  "I5" = 0.748;
    "J5" = 0.832;  "K5" = 0.787;

    "B6" = 0.975;  "C6" = 0.999;  "D6" = 0.982;  "E6" = 0.9429999999999999;
    "F6" = 0.468;  "G6" = 0.326;  "H6" = 0.384;  "I6" = 0.731;
    "J6" = 0.832;  "K6" = 0.778;

    "E7" = 0.9350000000000001;
    "F7" = 0.493;  "G7" = 0.382;  "H7" = 0.43;   "I7" = 0.745;
    "J7" = 0.821;  "K7" = 0.782;

    "B8" = 0.929;  "C8" = 0.977;  "D8" = 0.949;  "E8" = 0.833;
    "F8" = 0.41;   "G8" = 0.281;  "H8" = 0.333;  "I8" = 0.714;
    "J8" = 0.8159999999999999;  "K8" = 0.762;

    "B9" = 0.975;  "C9" = 0.999;  "D9" = 0.982;  "E9" = 0.945;
    "F9" = 0.462;  "G9" = 0.404;  "H9" = 0.431;  "I9" = 0.744;
    "J9" = 0.786;  "K9" = 0.764;
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
